$wb = $excel.ActiveWorkbook

# --- Text change: "Ready for handoff" -> "In Translation" ---
# This shared string is used as the Status value for the "Latest HO" /
# "Latest Target" columns on every sheet (Overview!E2:F2, zh-cn!C2,
# de-de!C2); updating the cell values re-points them to (or creates) the
# shared string "In Translation".
$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"
$wsZhCn.Range("C2").Value     = "In Translation"
$wsDeDe.Range("C2").Value     = "In Translation"

# --- Column width changes (Status columns got narrower after the copy
# edit): Overview!E:F and zh-cn!C / de-de!C shrink from ~17.22 to ~13.41
# "characters". ColumnWidth snaps to whole-pixel increments, so feed it
# the character value whose pixel rounding lands closest to the target
# stored width (13.4101845877511).
$wsOverview.Range("E1:F1").ColumnWidth = 12.5
$wsZhCn.Range("C1").ColumnWidth = 12.5
$wsDeDe.Range("C1").ColumnWidth = 12.5
